$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: fill in accuracy, date and logfile that had been lost ---
# (formulas in D/E/F already exist on the sheet and recalc automatically)
$ws.Range("C20").Value = 0.8125
$ws.Range("G20").Value = 42538
$ws.Range("H20").Value = "11.07.txt"

# --- Restore the three comments describing the "15." meta-group tests ---
$c = $ws.Range("B20").AddComment("15.Start`nPart of the 15. meta-group.`nTests the first 15 seconds of the sample against a specific neural network form. The other meta-group members test different parts of the song.`n")
$c = $ws.Range("B23").AddComment("15.Minute`nPart of the 15. meta-group.`nTests the 15 seconds from 1:00 to 1:15 of the sample against a specific neural network form. The other meta-group members test different parts of the song.`n")
$c = $ws.Range("B26").AddComment("15.Mix`nPart of the 15. meta-group.`nTests a composite 15 seconds (made up of 3 5-second chunks taken from random points within the sample) against a specific neural network form. The other meta-group members test different parts of the song.`n")

# --- Restore the selected cell ---
$ws.Range("G21").Select() | Out-Null
